# Re-sort the DEU electricity-region rows ("rez_DEU_*" / "elc*DEU_*") in the
# geo_sets worksheet from lexicographic/string order into numeric order
# (rez_DEU_0, rez_DEU_1, rez_DEU_2, ... rez_DEU_188) while keeping the pairing
# between column B ("rez_DEU_N") and column C ("elc*DEU_NNNN") intact.
#
# Note: Sort-Object with a scriptblock key selector is unreliable for
# double-digit+ numeric keys in this host, so the sort is implemented
# manually (simple selection sort) instead of relying on Sort-Object.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("geo_sets")

$firstRow = 642
$lastRow  = 829
$rowCount = $lastRow - $firstRow + 1

$rng = $ws.Range("B$firstRow" + ":C$lastRow")
$vals = $rng.Value()

# Pull the current contents into three parallel arrays: numeric sort key,
# column B text, column C text.
$keys = @()
$colB = @()
$colC = @()
for ($i = 1; $i -le $rowCount; $i++) {
    $b = [string]$vals[$i, 1]
    $c = [string]$vals[$i, 2]

    $m = [System.Text.RegularExpressions.Regex]::Match($b, '_(\d+)$')
    $key = [int]$m.Groups[1].Value

    $keys += $key
    $colB += $b
    $colC += $c
}

# Manual in-place selection sort (stable for our purposes since all keys
# here are unique) on the three parallel arrays, keyed on $keys.
for ($i = 0; $i -lt $rowCount - 1; $i++) {
    $minIdx = $i
    for ($j = $i + 1; $j -lt $rowCount; $j++) {
        if ($keys[$j] -lt $keys[$minIdx]) {
            $minIdx = $j
        }
    }
    if ($minIdx -ne $i) {
        $tmp = $keys[$i]; $keys[$i] = $keys[$minIdx]; $keys[$minIdx] = $tmp
        $tmp = $colB[$i]; $colB[$i] = $colB[$minIdx]; $colB[$minIdx] = $tmp
        $tmp = $colC[$i]; $colC[$i] = $colC[$minIdx]; $colC[$minIdx] = $tmp
    }
}

# Write the new, numerically-sorted order back into the same B:C range.
$out = New-Object 'object[,]' $rowCount,2
for ($i = 0; $i -lt $rowCount; $i++) {
    $out[$i, 0] = $colB[$i]
    $out[$i, 1] = $colC[$i]
}

$rng.Value = $out
